$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 17: drop formula, keep cached value only ---
$ws.Range("C17").Value = 113207000

# --- Row 19: drop formula, keep cached value only ---
$ws.Range("C19").Value = 8173500

# --- Row 20: Debit value changes 10000 -> 100000 ---
$ws.Range("D20").Value = 100000

# --- Row 22: Debit formula changes 45000 -> 45000+180000 ---
$ws.Range("D22").Formula = "=45000+180000"

# --- Row 23: Debit formula changes 450000+487500 -> 450000+487500+485000 ---
$ws.Range("D23").Formula = "=450000+487500+485000"

# --- Row 25: new transaction: GARRETH - lab ---
$ws.Range("B25").Value = "GARRETH - lab"
$ws.Range("D25").Formula = "=1150000"

# --- Row 26: new transaction: A/R ---
$ws.Range("B26").Value = "A/R"
$ws.Range("C26").Formula = "=9560000+18340500"

# --- Row 27: new transaction: SALES - cash/retail ---
$ws.Range("B27").Value = "SALES - cash/retail"
$ws.Range("C27").Formula = "=37514525-7353025-18340500"

# --- Row 28: new transaction: SELISIH - kurang ---
$ws.Range("B28").Value = "SELISIH - kurang"
$ws.Range("D28").Value = 3000

# --- Row 29: new transaction: SETOR KE BANK ---
$ws.Range("B29").Value = "SETOR KE BANK"
$ws.Range("D29").Value = 37000000

# --- Row 30: new day (5 Mar 2021), Wages Expense ---
$ws.Range("A30").Value = 44260
$ws.Range("B30").Value = "Wages Expense"
$ws.Range("D30").Formula = "=45000+195000"

# --- Row 31: new transaction: A/P ---
$ws.Range("B31").Value = "A/P"
$ws.Range("D31").Formula = "=2600000"

# --- Row 32: new transaction: TRANSFER BCA ---
$ws.Range("B32").Value = "TRANSFER BCA"
$ws.Range("D32").Formula = "=3700000+4500000+1649500+5020000+2100000"

# --- Row 33: new transaction: A/R ---
$ws.Range("B33").Value = "A/R"
$ws.Range("C33").Formula = "=29190500"

# --- Row 34: new transaction: SALES - cash/retail ---
$ws.Range("B34").Value = "SALES - cash/retail"
$ws.Range("C34").Formula = "=19282975+20997025-29190500"

# --- Row 35: new transaction: QIU - dokter ---
$ws.Range("B35").Value = "QIU - dokter"
$ws.Range("D35").Formula = "=300000"

# --- Row 36: new transaction: SELISIH - lebih ---
$ws.Range("B36").Value = "SELISIH - lebih"
$ws.Range("C36").Value = 11000

# --- Row 37: new transaction: SETOR KE BANK ---
$ws.Range("B37").Value = "SETOR KE BANK"
$ws.Range("D37").Value = 20000000

# --- Row 38: new day (6 Mar 2021), Wages Expense ---
$ws.Range("A38").Value = 44261
$ws.Range("B38").Value = "Wages Expense"
$ws.Range("D38").Formula = "=45000+1180000"

# --- Row 39: new transaction: BELI stempet ---
$ws.Range("B39").Value = "BELI stempet"
$ws.Range("D39").Formula = "=18000"

# --- Row 40: new transaction: TRANSFER BCA ---
$ws.Range("B40").Value = "TRANSFER BCA"
$ws.Range("D40").Formula = "=85000+200000+840000+1280000+2907000+30000000"

# --- Row 41: new transaction: IURAN DAERAH ---
$ws.Range("B41").Value = "IURAN DAERAH"
$ws.Range("D41").Value = 25000

# --- Row 42: new transaction: PARKIR - bulanan ---
$ws.Range("B42").Value = "PARKIR - bulanan"
$ws.Range("D42").Value = 10000

# --- Row 43: new transaction: BENSIN - RUSH ---
$ws.Range("B43").Value = "BENSIN - RUSH"
$ws.Range("D43").Value = 250000

# --- Row 44: new transaction: A/R ---
$ws.Range("B44").Value = "A/R"
$ws.Range("C44").Formula = "=47072000+30000000"

# --- Row 45: new transaction: SALES - cash/retail ---
$ws.Range("B45").Value = "SALES - cash/retail"
$ws.Range("C45").Formula = "=5146975+47526025-47072000"

# --- Row 46: new transaction: SELISIH - lebih ---
$ws.Range("B46").Value = "SELISIH - lebih"
$ws.Range("C46").Value = 140000

# --- Row 47: new transaction: SETOR KE BANK ---
$ws.Range("B47").Value = "SETOR KE BANK"
$ws.Range("D47").Value = 46000000

# --- Update selection to reflect scrolled view (D21 -> D41) ---
$ws.Activate()
$ws.Range("D41").Select()
